# Update cryptos list (GitHub Actions style refresh of coinranking data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to force text assignment (columns hold text, never real numbers,
# even when the value looks numeric, e.g. "680.21" or "7.10").
function Set-Text($row, $col, $value) {
    $text = [string]$value
    if ($text -match '^[0-9]+(\.[0-9]+)?$') {
        # Prefix with an apostrophe so Excel stores it as literal text
        # instead of silently converting it to a number.
        $ws.Cells.Item($row, $col).Value = "'" + $text
    } else {
        $ws.Cells.Item($row, $col).Value = $text
    }
}

# --- Row 49 / 50 swap: ONDO <-> SuiNetwork, plus new E values ---
Set-Text 49 2 "SuiNetwork"
Set-Text 49 3 "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-Text 49 4 "1.10"
Set-Text 49 5 "  -4.40%  "

Set-Text 50 2 "ONDO"
Set-Text 50 3 "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-Text 50 4 "1.29"
Set-Text 50 5 "  -5.18%  "

# --- D (Price) and E (Volume 1h) updates for all other rows ---
# Each entry: row, D new value (or $null to leave unchanged), E new value
$updates = @(
    @{ Row = 2;  D = "69.365.65";  E = "  -0.14%  " }
    @{ Row = 3;  D = "3.688.16";   E = "  -0.33%  " }
    @{ Row = 4;  D = $null;        E = "  +0.11%  " }
    @{ Row = 5;  D = "680.21";     E = "  -1.97%  " }
    @{ Row = 6;  D = "159.46";     E = "  -2.57%  " }
    @{ Row = 7;  D = $null;        E = "  +0.03%  " }
    @{ Row = 8;  D = $null;        E = "  -1.32%  " }
    @{ Row = 9;  D = $null;        E = "  -1.68%  " }
    @{ Row = 10; D = "7.10";       E = "  -4.22%  " }
    @{ Row = 11; D = $null;        E = "  -1.88%  " }
    @{ Row = 12; D = $null;        E = "  -3.56%  " }
    @{ Row = 13; D = "4.312.79";   E = "  -0.14%  " }
    @{ Row = 14; D = $null;        E = "  -3.18%  " }
    @{ Row = 15; D = "3.681.36";   E = "  -0.31%  " }
    @{ Row = 16; D = "69.357.47";  E = "  -0.21%  " }
    @{ Row = 17; D = $null;        E = "  +1.62%  " }
    @{ Row = 18; D = "16.03";      E = "  -1.85%  " }
    @{ Row = 19; D = "6.42";       E = "  -3.02%  " }
    @{ Row = 20; D = "468.69";     E = "  -2.91%  " }
    @{ Row = 21; D = "10.02";      E = "  +0.17%  " }
    @{ Row = 22; D = $null;        E = "  -2.33%  " }
    @{ Row = 23; D = "79.91";      E = "  -0.34%  " }
    @{ Row = 24; D = "3.835.35";   E = "  -0.13%  " }
    @{ Row = 25; D = $null;        E = "  -0.08%  " }
    @{ Row = 26; D = $null;        E = "  -6.39%  " }
    @{ Row = 27; D = $null;        E = "  -4.55%  " }
    @{ Row = 28; D = $null;        E = "  -4.50%  " }
    @{ Row = 29; D = $null;        E = "  -2.20%  " }
    @{ Row = 30; D = $null;        E = "  -4.16%  " }
    @{ Row = 31; D = $null;        E = "  -3.67%  " }
    @{ Row = 32; D = $null;        E = "  -4.69%  " }
    @{ Row = 33; D = "0.999";      E = "  +0.19%  " }
    @{ Row = 34; D = "26.93";      E = "  -0.85%  " }
    @{ Row = 35; D = "3.678.81";   E = "  +0.38%  " }
    @{ Row = 36; D = $null;        E = "  -5.30%  " }
    @{ Row = 37; D = $null;        E = "  -3.02%  " }
    @{ Row = 38; D = "6.23";       E = "  -2.53%  " }
    @{ Row = 39; D = $null;        E = "  -0.01%  " }
    @{ Row = 40; D = $null;        E = "  -2.88%  " }
    @{ Row = 41; D = $null;        E = "  -0.03%  " }
    @{ Row = 42; D = "0.0905";     E = "  -3.42%  " }
    @{ Row = 43; D = "171.69";     E = "  +4.74%  " }
    @{ Row = 44; D = "0.943";      E = "  -1.17%  " }
    @{ Row = 45; D = $null;        E = "  -0.74%  " }
    @{ Row = 46; D = "28.36";      E = "  -6.20%  " }
    @{ Row = 47; D = $null;        E = "  -4.55%  " }
    @{ Row = 48; D = $null;        E = "  -3.50%  " }
    @{ Row = 51; D = "7.79";       E = "  -3.08%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-Text $u.Row 4 $u.D
    }
    Set-Text $u.Row 5 $u.E
}
